$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 342; everything from old row 342 downward
# shifts down by one (old 342 -> new 343, ..., old 364 -> new 365).
$ws.Rows.Item(342).Insert()

# Populate the newly inserted row 342 with its data.
$ws.Range("A342").Value = 6
$ws.Range("B342").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C342").Value = "Metropolitana"
$ws.Range("D342").Value = 44585
$ws.Range("D342").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E342").Value = 13
$ws.Range("F342").Value = 100112039
$ws.Range("G342").Value = "Ciboulette"
$ws.Range("H342").Value = "Sin especificar"
$ws.Range("I342").Value = "Primera"
$ws.Range("J342").Value = 690
$ws.Range("K342").Value = 900
$ws.Range("L342").Value = 1000
$ws.Range("M342").Value = 948
$ws.Range("N342").Value = "`$/docena de atados"
$ws.Range("O342").Value = "Región Metropolitana"
$ws.Range("P342").Value = 316
$ws.Range("Q342").Value = 3
$ws.Range("R342").Value = "Hortaliza"
